$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$comment = "Valor Agregado Bruto (vab) en pesos corrientes y luego ajustado o convertido por distintas variables (IPC, IPI, TCP, TCC). Empleo en puestos de trabajo asalariados (emp). La productividad del trabajo se presenta tanto en índice como en nivel."

$lastRow = 253

# New header in X1, copying format/style from the existing header cell W1
$ws.Range("W1").Copy()
$ws.Range("X1").PasteSpecial(-4122)
$ws.Range("X1").Value = "unidades"

# Fill the comment text for every data row in column X (no special style, matches data cells)
$ws.Range("X2:X$lastRow").Value = $comment

$excel.CutCopyMode = 0
